# Refresh the cryptos snapshot: updated Price (D) and Volume(1h) (E) columns
# for the latest scrape. Values are assigned as literal text - a leading
# apostrophe is used for price strings that would otherwise be
# auto-coerced into numbers by Excel (dropping trailing zeros / exponent
# notation), matching the inline-string storage of the source cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.199.13"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.602.55"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'303.00"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").Value = "'0.3783"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "'51.96"
$ws.Range("E8").Value = "  +3.93%  "
$ws.Range("D9").Value = "'0.3615"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").Value = "'1.266"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "'0.08134"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'22.64"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "'6.584"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "'7.391"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "'0.00001249"
$ws.Range("D17").Value = "1.600.30"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "'93.86"
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'18.07"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("D21").Value = "'6.548"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "23.190.21"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'2.388"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("D26").Value = "'2.983"
$ws.Range("E26").Value = "  +9.90%  "
$ws.Range("D27").Value = "'21.21"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").Value = "'148.90"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").Value = "'133.83"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "'2.375"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").Value = "'6.837"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "1.778.17"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value = "'0.9715"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").Value = "'10.34"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("D37").Value = "'0.02716"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'0.2510"
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("D39").Value = "'6.119"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").Value = "'0.08809"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "'1.363"
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Value = "'0.7096"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("D43").Value = "'12.53"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "'15.55"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").Value = "'0.6535"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "'2.312"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "'1.218"
$ws.Range("E51").Value = "  +3.55%  "
